$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns D (Price) and E (Volume) keep their exact text representation
# (some values like "6.65" or "0.0000118" would otherwise be auto-converted to
# numbers by Excel, losing trailing zeros / formatting).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 7 and Row 8 swap contents (XRP <-> USDC)
$ws.Range("B7").Value = "USDC"
$ws.Range("C7").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("B8").Value = "XRP"
$ws.Range("C8").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D8").Value = "0.603"
$ws.Range("E8").Value = "  +0.68%  "

$ws.Range("D2").Value = "67.188.25"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "3.316.78"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "184.78"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "576.79"
$ws.Range("E9").Value = "  -0.20%  "
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  +0.68%  "
$ws.Range("D11").Value = "0.406"
$ws.Range("E11").Value = "  -0.03%  "
$ws.Range("D12").Value = "3.899.82"
$ws.Range("E12").Value = "  +1.87%  "
$ws.Range("E13").Value = "  -0.79%  "
$ws.Range("D14").Value = "27.29"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("D15").Value = "67.345.95"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  -0.21%  "
$ws.Range("D17").Value = "3.318.77"
$ws.Range("E17").Value = "  +1.73%  "
$ws.Range("D18").Value = "442.84"
$ws.Range("E18").Value = "  +6.71%  "
$ws.Range("D19").Value = "13.56"
$ws.Range("E19").Value = "  +2.43%  "
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").Value = "7.70"
$ws.Range("E21").Value = "  +2.51%  "
$ws.Range("D22").Value = "74.15"
$ws.Range("E22").Value = "  +4.17%  "
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "3.465.21"
$ws.Range("E24").Value = "  +1.84%  "
$ws.Range("D25").Value = "0.512"
$ws.Range("E25").Value = "  +1.19%  "
$ws.Range("D26").Value = "0.0000118"
$ws.Range("E26").Value = "  +1.58%  "
$ws.Range("E27").Value = "  +1.62%  "
$ws.Range("D28").Value = "9.01"
$ws.Range("E28").Value = "  -3.86%  "
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("E30").Value = "  +1.40%  "
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").Value = "5.32"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("E33").Value = "  -0.04%  "
$ws.Range("D34").Value = "6.81"
$ws.Range("E34").Value = "  -0.18%  "
$ws.Range("D35").Value = "1.23"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("D36").Value = "1.50"
$ws.Range("E36").Value = "  +4.54%  "
$ws.Range("D37").Value = "162.31"
$ws.Range("E37").Value = "  -0.34%  "
$ws.Range("D38").Value = "27.20"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("D40").Value = "2.789.77"
$ws.Range("E40").Value = "  +5.96%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").Value = "4.47"
$ws.Range("E42").Value = "  +0.64%  "
$ws.Range("E43").Value = "  -1.01%  "
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("E45").Value = "  -0.40%  "
$ws.Range("D46").Value = "24.62"
$ws.Range("E46").Value = "  +1.84%  "
$ws.Range("D47").Value = "2.38"
$ws.Range("E47").Value = "  -1.76%  "
$ws.Range("D48").Value = "325.73"
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("D49").Value = "0.0272"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("D50").Value = "0.988"
$ws.Range("E50").Value = "  +1.31%  "
$ws.Range("D51").Value = "31.08"
$ws.Range("E51").Value = "  +1.72%  "
